$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 182, shifting existing rows 182:270 down to 183:271
$ws.Rows(182).Insert()

# Populate the newly inserted row 182 with the new record's data.
# Columns A,B,C,E,F,G,H,I,R are constant for every data row in this sheet.
$ws.Cells.Item(182, 1).Value = 4
$ws.Cells.Item(182, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(182, 3).Value = "Los Lagos"
$ws.Cells.Item(182, 4).Value = 44960
$ws.Cells.Item(182, 5).Value = 10
$ws.Cells.Item(182, 6).Value = 100112009
$ws.Cells.Item(182, 7).Value = "Acelga"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 80
$ws.Cells.Item(182, 11).Value = 10000
$ws.Cells.Item(182, 12).Value = 10000
$ws.Cells.Item(182, 13).Value = 10000
$ws.Cells.Item(182, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(182, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(182, 16).Value = 833
$ws.Cells.Item(182, 17).Value = 12
$ws.Cells.Item(182, 18).Value = "Hortaliza"

# Match the date-formatted style (s="2") used by column D in every other row
$ws.Cells.Item(182, 4).NumberFormat = $ws.Cells.Item(183, 4).NumberFormat
